# Update the Gdf1-Bmpr2 LR-pairs sheet with refreshed TPM-derived values.
# Rows 2-6 (FAPs as sender) keep their sender/ligand/receptor/target labels but
# get refreshed numeric values; rows 7-11 are new (MuSCs as sender) mirroring the
# same five target clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Gdf1"
$ws.Cells.Item(2, 3).Value = "Bmpr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.1814253333333333
$ws.Cells.Item(2, 8).Value = 0.544276
$ws.Cells.Item(2, 9).Value = 0.9591577789839493
$ws.Cells.Item(2, 10).Value = 0.9591577789839494
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 45.924193
$ws.Cells.Item(2, 14).Value = 137.772579
$ws.Cells.Item(2, 15).Value = 0.307792367338991
$ws.Cells.Item(2, 16).Value = 0.307792367338991
$ws.Cells.Item(2, 17).Value = 8.331812023089334
$ws.Cells.Item(2, 18).Value = 74.986308207804
$ws.Cells.Item(2, 19).Value = 0.2952214434450784
$ws.Cells.Item(2, 20).Value = 0.2952214434450785

# Row 3: FAPs -> FAPs
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Gdf1"
$ws.Cells.Item(3, 3).Value = "Bmpr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.1814253333333333
$ws.Cells.Item(3, 8).Value = 0.544276
$ws.Cells.Item(3, 9).Value = 0.9591577789839493
$ws.Cells.Item(3, 10).Value = 0.9591577789839494
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 40.23702866666667
$ws.Cells.Item(3, 14).Value = 120.711086
$ws.Cells.Item(3, 15).Value = 0.2696759485354523
$ws.Cells.Item(3, 16).Value = 0.2696759485354523
$ws.Cells.Item(3, 17).Value = 7.300016338192889
$ws.Cells.Item(3, 18).Value = 65.700147043736
$ws.Cells.Item(3, 19).Value = 0.2586617838426542
$ws.Cells.Item(3, 20).Value = 0.2586617838426543

# Row 4: FAPs -> Inflammatory-Mac
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Gdf1"
$ws.Cells.Item(4, 3).Value = "Bmpr2"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.1814253333333333
$ws.Cells.Item(4, 8).Value = 0.544276
$ws.Cells.Item(4, 9).Value = 0.9591577789839493
$ws.Cells.Item(4, 10).Value = 0.9591577789839494
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 23.60320766666667
$ws.Cells.Item(4, 14).Value = 70.809623
$ws.Cells.Item(4, 15).Value = 0.1581930283351338
$ws.Cells.Item(4, 16).Value = 0.1581930283351339
$ws.Cells.Item(4, 17).Value = 4.282219818660889
$ws.Cells.Item(4, 18).Value = 38.539978367948
$ws.Cells.Item(4, 19).Value = 0.1517320737086719
$ws.Cells.Item(4, 20).Value = 0.1517320737086719

# Row 5: FAPs -> MuSCs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gdf1"
$ws.Cells.Item(5, 3).Value = "Bmpr2"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.1814253333333333
$ws.Cells.Item(5, 8).Value = 0.544276
$ws.Cells.Item(5, 9).Value = 0.9591577789839493
$ws.Cells.Item(5, 10).Value = 0.9591577789839494
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 19.226538
$ws.Cells.Item(5, 14).Value = 57.679614
$ws.Cells.Item(5, 15).Value = 0.1288597852280838
$ws.Cells.Item(5, 16).Value = 0.1288597852280838
$ws.Cells.Item(5, 17).Value = 3.488181065496
$ws.Cells.Item(5, 18).Value = 31.393629589464
$ws.Cells.Item(5, 19).Value = 0.1235968653997176
$ws.Cells.Item(5, 20).Value = 0.1235968653997176

# Row 6: FAPs -> Resolving-Mac
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gdf1"
$ws.Cells.Item(6, 3).Value = "Bmpr2"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.1814253333333333
$ws.Cells.Item(6, 8).Value = 0.544276
$ws.Cells.Item(6, 9).Value = 0.9591577789839493
$ws.Cells.Item(6, 10).Value = 0.9591577789839494
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 20.21413933333333
$ws.Cells.Item(6, 14).Value = 60.64241799999999
$ws.Cells.Item(6, 15).Value = 0.1354788705623391
$ws.Cells.Item(6, 16).Value = 0.1354788705623391
$ws.Cells.Item(6, 17).Value = 3.667356966596444
$ws.Cells.Item(6, 18).Value = 33.00621269936799
$ws.Cells.Item(6, 19).Value = 0.1299456125878271
$ws.Cells.Item(6, 20).Value = 0.1299456125878271

# Row 7: MuSCs -> ECs
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Gdf1"
$ws.Cells.Item(7, 3).Value = "Bmpr2"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.007725333333333334
$ws.Cells.Item(7, 8).Value = 0.023176
$ws.Cells.Item(7, 9).Value = 0.0408422210160507
$ws.Cells.Item(7, 10).Value = 0.0408422210160507
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 45.924193
$ws.Cells.Item(7, 14).Value = 137.772579
$ws.Cells.Item(7, 15).Value = 0.307792367338991
$ws.Cells.Item(7, 16).Value = 0.307792367338991
$ws.Cells.Item(7, 17).Value = 0.3547796989893334
$ws.Cells.Item(7, 18).Value = 3.193017290904001
$ws.Cells.Item(7, 19).Value = 0.01257092389391254
$ws.Cells.Item(7, 20).Value = 0.01257092389391254

# Row 8: MuSCs -> FAPs
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Gdf1"
$ws.Cells.Item(8, 3).Value = "Bmpr2"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.007725333333333334
$ws.Cells.Item(8, 8).Value = 0.023176
$ws.Cells.Item(8, 9).Value = 0.0408422210160507
$ws.Cells.Item(8, 10).Value = 0.0408422210160507
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 40.23702866666667
$ws.Cells.Item(8, 14).Value = 120.711086
$ws.Cells.Item(8, 15).Value = 0.2696759485354523
$ws.Cells.Item(8, 16).Value = 0.2696759485354523
$ws.Cells.Item(8, 17).Value = 0.3108444587928889
$ws.Cells.Item(8, 18).Value = 2.797600129136001
$ws.Cells.Item(8, 19).Value = 0.01101416469279806
$ws.Cells.Item(8, 20).Value = 0.01101416469279806

# Row 9: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Gdf1"
$ws.Cells.Item(9, 3).Value = "Bmpr2"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.007725333333333334
$ws.Cells.Item(9, 8).Value = 0.023176
$ws.Cells.Item(9, 9).Value = 0.0408422210160507
$ws.Cells.Item(9, 10).Value = 0.0408422210160507
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 23.60320766666667
$ws.Cells.Item(9, 14).Value = 70.809623
$ws.Cells.Item(9, 15).Value = 0.1581930283351338
$ws.Cells.Item(9, 16).Value = 0.1581930283351339
$ws.Cells.Item(9, 17).Value = 0.1823426469608889
$ws.Cells.Item(9, 18).Value = 1.641083822648
$ws.Cells.Item(9, 19).Value = 0.006460954626461906
$ws.Cells.Item(9, 20).Value = 0.006460954626461908

# Row 10: MuSCs -> MuSCs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Gdf1"
$ws.Cells.Item(10, 3).Value = "Bmpr2"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.007725333333333334
$ws.Cells.Item(10, 8).Value = 0.023176
$ws.Cells.Item(10, 9).Value = 0.0408422210160507
$ws.Cells.Item(10, 10).Value = 0.0408422210160507
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 19.226538
$ws.Cells.Item(10, 14).Value = 57.679614
$ws.Cells.Item(10, 15).Value = 0.1288597852280838
$ws.Cells.Item(10, 16).Value = 0.1288597852280838
$ws.Cells.Item(10, 17).Value = 0.148531414896
$ws.Cells.Item(10, 18).Value = 1.336782734064
$ws.Cells.Item(10, 19).Value = 0.005262919828366224
$ws.Cells.Item(10, 20).Value = 0.005262919828366224

# Row 11: MuSCs -> Resolving-Mac
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Gdf1"
$ws.Cells.Item(11, 3).Value = "Bmpr2"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.007725333333333334
$ws.Cells.Item(11, 8).Value = 0.023176
$ws.Cells.Item(11, 9).Value = 0.0408422210160507
$ws.Cells.Item(11, 10).Value = 0.0408422210160507
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 20.21413933333333
$ws.Cells.Item(11, 14).Value = 60.64241799999999
$ws.Cells.Item(11, 15).Value = 0.1354788705623391
$ws.Cells.Item(11, 16).Value = 0.1354788705623391
$ws.Cells.Item(11, 17).Value = 0.1561609643964444
$ws.Cells.Item(11, 18).Value = 1.405448679568
$ws.Cells.Item(11, 19).Value = 0.005533257974511979
$ws.Cells.Item(11, 20).Value = 0.005533257974511979
